$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $null) { continue }
    $text = [string]$val
    if ($text.IndexOf(",") -lt 0) { continue }

    $parts = $text.Split(",")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $sysIndex = -1
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i].Equals("System")) {
            $sysIndex = $i
            break
        }
    }

    if ($sysIndex -ge 0) {
        $rest = @()
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $sysIndex) { $rest += $parts[$i] }
        }
        $newParts = @("System") + $rest
    } else {
        $newParts = @()
        for ($i = $parts.Length - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $newText = [string]::Join(", ", $newParts)
    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}

Write-Host "Done"
